# "Home pagge , usermanagement and order"
#
# 1) Bump the cached master/layout date field (datetimeFigureOut) from
#    8/16/2021 -> 8/17/2021 everywhere it appears (slide master + all
#    custom layouts).
# 2) Slide 1: rename/move the "Khmershop" textbox to "T-shop" at its new
#    position.
# 3) Slide 1: remove the standalone "LOGO" textbox.

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text on the master and every layout ------------
function Update-DatePlaceholder($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "8/17/2021"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($i)
}

# --- 2) & 3) Slide 1 shape tweaks ----------------------------------------
$slide1 = $p.Slides.Item(1)

for ($j = 1; $j -le $slide1.Shapes.Count; $j++) {
    $shp = $slide1.Shapes.Item($j)

    if ($shp.Name -eq "TextBox 4") {
        # Khmershop -> T-shop, shifted up and to the left.
        $shp.Left = 35.79236220472441
        $shp.Top = 8.304409548818898
        $shp.TextFrame.TextRange.Text = "T-shop"
    }
}

# Remove the leftover "LOGO" textbox (walk backwards since we delete).
for ($j = $slide1.Shapes.Count; $j -ge 1; $j--) {
    $shp = $slide1.Shapes.Item($j)
    if ($shp.Name -eq "TextBox 25") {
        $shp.Delete()
    }
}
